$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "DATA201"
$ws.Range("B3").Value = "Intro to Databases"
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = "Utilizing SQL in Data Science"

$ws.Range("D3").Select()
